# LV Contacts Regression Suite - 21 June 2024
#
# On the "Users" sheet, add a second user (Admin User / Indrajeet Singh)
# next to the existing one, matching the header's bold style, and leave
# the "Users" sheet as the active/selected tab (moving it off the
# "RecentlyViewedListView" sheet, which was previously active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New header cell B1 "Admin User" - bold, like A1 "Users"
$ws.Range("B1").Value = "Admin User"
$ws.Range("B1").Font.Bold = $true

# New data cell B2 "Indrajeet Singh" under the new header
$ws.Range("B2").Value = "Indrajeet Singh"

# Make "Users" the active sheet/tab, with D8 selected
$ws.Activate()
$ws.Range("D8").Select()
